# Daily attendance processing - 2025-12-25 19:51:45
#
# Reorders the "Recorded By" list in column G of the Session Analysis
# Results sheet: any entries that are literally "System" (any letter
# case) are moved to the end of the comma-separated list, while the
# relative order of all entries (system and non-system alike) is
# otherwise preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -eq "") { continue }

    $parts = $val -split ", "

    $nonSystemParts = @()
    $systemParts = @()
    foreach ($part in $parts) {
        if ($part.ToLower() -eq "system") {
            $systemParts += $part
        } else {
            $nonSystemParts += $part
        }
    }

    if ($systemParts.Count -eq 0) { continue }

    $reordered = $nonSystemParts + $systemParts
    $newVal = [string]::Join(", ", $reordered)

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
